$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "38.263.24"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +3.30%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.122.92"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.71%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "235.44"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.78%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.627"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.65%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "58.34"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.29%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.393"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.50%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0783"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.97%  "
$ws.Range("E11").Value = "  +1.82%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "2.436.56"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.70%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "14.61"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.75%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "21.65"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.793"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.13%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.48%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.109.42"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.90%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "38.204.54"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +3.35%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.19"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.07%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "70.62"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.84%  "
$ws.Range("E21").Value = "  +2.97%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "229.43"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.41%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.41"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.48%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.41"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.26%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "169.15"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("E27").Value = "  +11.64%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "9.02"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.35%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.43"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.57%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "19.61"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +3.60%  "
$ws.Range("E31").Value = "  +1.10%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.65"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.00%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.63"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +4.57%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0627"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.86%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "4.61"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("E36").Value = "  +7.89%  "
$ws.Range("E37").Value = "  +4.98%  "
$ws.Range("E38").Value = "  -0.06%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.49"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.76%  "
$ws.Range("E40").Value = "  +8.99%  "
$ws.Range("E41").Value = "  +0.21%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "97.43"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.70%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0216"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +3.20%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.461.99"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  +5.73%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "4.12"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -8.39%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "15.76"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.95%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "3.05"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.05%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "7.33"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.34%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.321.10"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.69%  "
